# Apply updated cryptocurrency price/volume data to the active worksheet.
# Column D ("Price") values that look like plain numbers are written with a
# leading apostrophe so Excel stores them as text (matching the source data,
# which keeps values such as "1.00" or "6.18" as text, not numbers) without
# altering the cells number format/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.442.65'
$ws.Range("E2").Value = '  +8.17%  '
# Row 3
$ws.Range("D3").Value = '2.578.88'
$ws.Range("E3").Value = '  +9.77%  '
# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.01%  '
# Row 5
$ws.Range("D5").Value = '''504.72'
$ws.Range("E5").Value = '  +6.95%  '
# Row 6
$ws.Range("D6").Value = '''156.79'
$ws.Range("E6").Value = '  +8.64%  '
# Row 7
$ws.Range("D7").Value = '''0.634'
$ws.Range("E7").Value = '  +26.35%  '
# Row 8
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.14%  '
# Row 9
$ws.Range("D9").Value = '2.574.75'
$ws.Range("E9").Value = '  +9.81%  '
# Row 10
$ws.Range("D10").Value = '''6.18'
$ws.Range("E10").Value = '  +15.06%  '
# Row 11
$ws.Range("D11").Value = '''0.103'
$ws.Range("E11").Value = '  +6.70%  '
# Row 12
$ws.Range("E12").Value = '  +6.79%  '
# Row 13
$ws.Range("D13").Value = '''0.127'
$ws.Range("E13").Value = '  +1.95%  '
# Row 14
$ws.Range("D14").Value = '3.023.81'
$ws.Range("E14").Value = '  +9.65%  '
# Row 15
$ws.Range("D15").Value = '59.313.16'
$ws.Range("E15").Value = '  +7.97%  '
# Row 16
$ws.Range("D16").Value = '''21.74'
$ws.Range("E16").Value = '  +8.81%  '
# Row 17
$ws.Range("E17").Value = '  +5.54%  '
# Row 18
$ws.Range("D18").Value = '2.575.86'
$ws.Range("E18").Value = '  +9.61%  '
# Row 19
$ws.Range("E19").Value = '  +5.48%  '
# Row 20
$ws.Range("D20").Value = '''336.19'
$ws.Range("E20").Value = '  +7.92%  '
# Row 21
$ws.Range("E21").Value = '  +7.75%  '
# Row 22
$ws.Range("D22").Value = '''6.06'
$ws.Range("E22").Value = '  +8.52%  '
# Row 23
$ws.Range("E23").Value = '  +0.76%  '
# Row 24
$ws.Range("D24").Value = '''59.95'
$ws.Range("E24").Value = '  +7.53%  '
# Row 25
$ws.Range("D25").Value = '''0.415'
$ws.Range("E25").Value = '  +6.19%  '
# Row 26
$ws.Range("E26").Value = '  +9.23%  '
# Row 27
$ws.Range("D27").Value = '2.694.88'
$ws.Range("E27").Value = '  +10.06%  '
# Row 28
$ws.Range("E28").Value = '  +0.18%  '
# Row 29
$ws.Range("D29").Value = '0.0₃0825'
$ws.Range("E29").Value = '  +9.90%  '
# Row 30
$ws.Range("E30").Value = '  +2.99%  '
# Row 31
$ws.Range("E31").Value = '  +0.06%  '
# Row 32
$ws.Range("D32").Value = '''157.80'
$ws.Range("E32").Value = '  +6.35%  '
# Row 33
$ws.Range("D33").Value = '''19.26'
$ws.Range("E33").Value = '  +7.37%  '
# Row 34
$ws.Range("D34").Value = '''1.57'
$ws.Range("E34").Value = '  +7.08%  '
# Row 35
$ws.Range("E35").Value = '  +9.65%  '
# Row 36
$ws.Range("E36").Value = '  +10.80%  '
# Row 37
$ws.Range("E37").Value = '  +9.45%  '
# Row 38
$ws.Range("D38").Value = '''0.849'
$ws.Range("E38").Value = '  +3.36%  '
# Row 39
$ws.Range("E39").Value = '  +12.05%  '
# Row 40
$ws.Range("D40").Value = '''1.45'
$ws.Range("E40").Value = '  +8.92%  '
# Row 41
$ws.Range("D41").Value = '''35.06'
$ws.Range("E41").Value = '  +4.77%  '
# Row 42
$ws.Range("D42").Value = '''291.78'
$ws.Range("E42").Value = '  +14.45%  '
# Row 43
$ws.Range("E43").Value = '  +8.48%  '
# Row 44
$ws.Range("D44").Value = '''0.621'
$ws.Range("E44").Value = '  +8.65%  '
# Row 45
$ws.Range("E45").Value = '  +7.33%  '
# Row 46
$ws.Range("E46").Value = '  +0.05%  '
# Row 47
$ws.Range("B47").Value = 'SuiNetwork'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D47").Value = '''0.743'
$ws.Range("E47").Value = '  +18.42%  '
# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''19.27'
$ws.Range("E48").Value = '  +15.00%  '
# Row 49
$ws.Range("D49").Value = '''0.0235'
# Row 50
$ws.Range("E50").Value = '  +7.62%  '
# Row 51
$ws.Range("D51").Value = '''10.25'
$ws.Range("E51").Value = '  +1.09%  '
